$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.734.07"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "1.964.78"
$ws.Range("E3").Value = "  +1.13%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.40%  "

$ws.Range("E6").Value = "  +0.83%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "59.20"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.08%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0820"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.90%  "

$ws.Range("E11").Value = "  -0.21%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.50"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.52%  "

$ws.Range("D13").Value = "2.254.93"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.76"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.19%  "

$ws.Range("E16").Value = "  +0.58%  "

$ws.Range("D17").Value = "1.966.22"
$ws.Range("E17").Value = "  +1.30%  "

$ws.Range("D18").Value = "36.628.24"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("E20").Value = "  -0.59%  "

$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "229.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "

$ws.Range("E23").Value = "  -0.10%  "

$ws.Range("E24").Value = "  -0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.03%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.144"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.91%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.43%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "160.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.91%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.119"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.52%  "

$ws.Range("E31").Value = "  -0.96%  "

$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0622"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.01%  "

$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("E35").Value = "  +6.69%  "

$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("E37").Value = "  -2.75%  "

$ws.Range("E38").Value = "  +12.49%  "

$ws.Range("E39").Value = "  -0.46%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.90"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.02%  "

$ws.Range("E42").Value = "  -0.77%  "

$ws.Range("E43").Value = "  +1.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.23"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.68%  "

$ws.Range("D45").Value = "1.362.00"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("E46").Value = "  +0.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.65%  "

$ws.Range("D50").Value = "2.145.94"
$ws.Range("E50").Value = "  +1.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "44.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.94%  "

